$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 44326
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44280
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("S6").Value = 900

# Row 7
$ws.Range("D7").Value = 44270

# Row 8
$ws.Range("D8").Value = 44364
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 23000
$ws.Range("P8").Value = 23000
$ws.Range("S8").Value = 1150

# Row 9
$ws.Range("D9").Value = 44306
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 1000

# Row 10
$ws.Range("D10").Value = 44285
$ws.Range("M10").Value = 70

# Row 11
$ws.Range("D11").Value = 44445
$ws.Range("M11").Value = 170
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("S11").Value = 1025

# Row 12
$ws.Range("D12").Value = 44425
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("S12").Value = 975

# Row 13
$ws.Range("D13").Value = 44333
$ws.Range("M13").Value = 70
$ws.Range("N13").Value = 18000
$ws.Range("P13").Value = 19000
$ws.Range("S13").Value = 950

# Row 14
$ws.Range("D14").Value = 44453
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 21000
$ws.Range("P14").Value = 20500
$ws.Range("R14").Value = 'Ecuador'
$ws.Range("S14").Value = 1025

# Row 15
$ws.Range("D15").Value = 44453
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("S15").Value = 1000

# Row 16
$ws.Range("D16").Value = 44202
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 17000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 17500
$ws.Range("S16").Value = 875

# Row 17
$ws.Range("D17").Value = 44405
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 22667
$ws.Range("S17").Value = 1133

# Row 18
$ws.Range("D18").Value = 44308
$ws.Range("N18").Value = 19000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19500
$ws.Range("S18").Value = 975

# Row 19
$ws.Range("D19").Value = 44271
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 16000
$ws.Range("O19").Value = 17000
$ws.Range("P19").Value = 16500
$ws.Range("S19").Value = 825

# Row 20
$ws.Range("D20").Value = 44431
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20500
$ws.Range("R20").Value = 'Perú'
$ws.Range("S20").Value = 1025

# Row 21
$ws.Range("D21").Value = 44355
$ws.Range("M21").Value = 80
$ws.Range("R21").Value = 'Costa Rica'

# Row 22
$ws.Range("D22").Value = 44410
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = 21000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 21500
$ws.Range("S22").Value = 1075

# Row 23
$ws.Range("D23").Value = 44284
$ws.Range("M23").Value = 70
$ws.Range("N23").Value = 18000
$ws.Range("O23").Value = 19000
$ws.Range("P23").Value = 18500
$ws.Range("S23").Value = 925

# Row 24
$ws.Range("D24").Value = 44406
$ws.Range("M24").Value = 150
$ws.Range("N24").Value = 20000
$ws.Range("O24").Value = 21000
$ws.Range("P24").Value = 20500
$ws.Range("S24").Value = 1025

# Row 25
$ws.Range("D25").Value = 44334
$ws.Range("M25").Value = 250

# Row 26
$ws.Range("D26").Value = 44299
$ws.Range("M26").Value = 120
$ws.Range("N26").Value = 19000
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 19500
$ws.Range("S26").Value = 975

# Row 27
$ws.Range("D27").Value = 44340
$ws.Range("R27").Value = 'Perú'

# Row 28
$ws.Range("D28").Value = 44371
$ws.Range("M28").Value = 150
$ws.Range("R28").Value = 'Costa Rica'

# Row 29
$ws.Range("D29").Value = 44419
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 21000
$ws.Range("P29").Value = 20500
$ws.Range("S29").Value = 1025

# Row 30
$ws.Range("D30").Value = 44300
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 18000
$ws.Range("O30").Value = 18000
$ws.Range("P30").Value = 18000
$ws.Range("S30").Value = 900
